# Resume edit: split the job-title run to add "and computer programmer",
# and reorder/rewrap the "interested in" runs to lead with "reproducible
# research, data science" instead of "clinical data science, reproducible
# science".
#
# NOTE: plain Range.Text / Find.Execute replacement causes this host to
# coalesce any newly-written text into the neighboring run whenever the
# resolved run formatting (rPr) is identical — exactly mirroring the
# "merge adjacent like-formatted runs" behavior real Word applies on
# save. Range.InsertXML, however, splices literal <w:r> elements in
# without that coalescing pass, so it is used here to reproduce the
# fine-grained run split shown in the target diff. InsertXML always
# lands its inserted content at the end of the *receiving* range, so
# each target range below is stretched out to the end of its paragraph
# and the untouched paragraph tail is re-emitted verbatim as trailing
# runs inside the same InsertXML payload.

$d = $word.ActiveDocument

function Find-ParagraphContaining($doc, $needle) {
  for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
    $para = $doc.Paragraphs($i)
    if ($para.Range.Text.Contains($needle)) {
      return $para
    }
  }
  throw "no paragraph containing '$needle' was found"
}

function Esc($s) {
  $s.Replace('&','&amp;').Replace('<','&lt;').Replace('>','&gt;')
}

function RunXml($text, $withRFonts) {
  $t = Esc($text)
  $space = ''
  if ($text -match '^\s' -or $text -match '\s$' -or $text -eq '') {
    $space = ' xml:space="preserve"'
  }
  if ($withRFonts) {
    return '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t' + $space + '>' + $t + '</w:t></w:r>'
  } else {
    return '<w:r><w:t' + $space + '>' + $t + '</w:t></w:r>'
  }
}

function Package($bodyXml) {
  return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + `
    $bodyXml + `
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# --- Hunk 1: "General practitioner and researcher at ..." paragraph ---
$p3 = Find-ParagraphContaining $d "eneral practitioner and researcher at Shiraz University of Medical Sciences, Shiraz, Iran"
$p3Start = $p3.Range.Start
$p3End = $p3.Range.End
$full3 = $p3.Range.Text
$idx3 = $full3.IndexOf("eneral practitioner and researcher at Shiraz University of Medical Sciences, Shiraz, Iran")
if ($idx3 -lt 0) { throw "hunk1 anchor text not found" }
$r3Start = $p3Start + $idx3
$rng3 = $d.Range($r3Start, $p3End)

$body3 = `
  (RunXml "eneral practitioner" $true) + `
  (RunXml ", " $true) + `
  (RunXml "researcher" $true) + `
  (RunXml ", and computer programmer" $true) + `
  (RunXml " at Shiraz University of Medical Sciences, Shiraz, Iran" $true) + `
  (RunXml "." $true)

$rng3.InsertXML( (Package $body3) )

# --- Hunk 2: "Extremely interested in ..." paragraph ---
$p4 = Find-ParagraphContaining $d "clinical data science, reproducible science,"
$p4Start = $p4.Range.Start
$p4End = $p4.Range.End
$full4 = $p4.Range.Text
$idx4 = $full4.IndexOf("clinical data science, reproducible science,")
if ($idx4 -lt 0) { throw "hunk2 anchor text not found" }
$r4Start = $p4Start + $idx4
$rng4 = $d.Range($r4Start, $p4End)

$body4 = `
  (RunXml "reproducible " $true) + `
  (RunXml "research," $true) + `
  (RunXml " " $false) + `
  (RunXml "data science" $false) + `
  (RunXml "," $true) + `
  (RunXml " clinical research, epidemiology, research design, and computer programming" $true) + `
  (RunXml "." $true)

$rng4.InsertXML( (Package $body4) )

Write-Output "P3 AFTER: $($p3.Range.Text)"
Write-Output "P4 AFTER: $($p4.Range.Text)"
